# The sheet gained five new sending/target-cluster combinations for the
# MuSCs <-> Resolving-Mac pair (Proc-Tek ligand/receptor), and the TPM-derived
# metrics for the existing rows were recomputed against the new totals.
# Rebuild the whole A2:T9 data block with the updated values in one pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# One inner array per data row (A..T), in sheet-row order starting at row 2.
$rows = @(
    @("MuSCs", "Proc", "Tek", "ECs", 2, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 7.44821386452, 67.03392478068, 0.6226736846770882, 0.6226736846770882),
    @("MuSCs", "Proc", "Tek", "FAPs", 2, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 1.403297318885, 12.629675869965, 0.1173162221361529, 0.1173162221361528),
    @("MuSCs", "Proc", "Tek", "MuSCs", 2, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 0.4326632403866666, 3.89396916348, 0.03617082149040259, 0.03617082149040259),
    @("MuSCs", "Proc", "Tek", "Resolving-Mac", 2, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.001250686158333333, 0.011256175425, 0.0001045578675303299, 0.0001045578675303299),
    @("Resolving-Mac", "Proc", "Tek", "ECs", 1, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 2.146719719664, 19.320477476976, 0.1794666348370614, 0.1794666348370614),
    @("Resolving-Mac", "Proc", "Tek", "FAPs", 1, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 0.4044575091153333, 3.640117582038, 0.03381281097416391, 0.0338128109741639),
    @("Resolving-Mac", "Proc", "Tek", "MuSCs", 1, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 0.1247019388817778, 1.122317449936, 0.01042513241191654, 0.01042513241191654),
    @("Resolving-Mac", "Proc", "Tek", "Resolving-Mac", 1, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.0003604720122222222, 0.00324424811, 0.0000301356056842818, 0.00003013560568428181)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
  $rowVals = $rows[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowVals[$c]
  }
}
